$d = $word.ActiveDocument

function Get-EscapedXml($text) {
    $t = $text -replace "&", "&amp;"
    $t = $t -replace "<", "&lt;"
    $t = $t -replace ">", "&gt;"
    return $t
}

function New-RunParagraphXml($text) {
    $escaped = Get-EscapedXml $text
    return '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr/><w:r><w:rPr><w:rFonts w:ascii="Helvetica" w:hAnsi="Helvetica" w:cs="Helvetica"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr><w:t xml:space="preserve">' + $escaped + '</w:t></w:r></w:p>'
}

# The final paragraph of the document currently reads " a. " (under
# "3. Identify Potential solutions" in the Predicting Fingers section).
# Replace its text with the completed answer, scoping the Find to that
# paragraph only so earlier " a. " occurrences are untouched.
$count = $d.Paragraphs.Count
$targetPara = $d.Paragraphs.Item($count)
$targetRange = $targetPara.Range
$targetRange.Find.Execute(" a. ", $false, $false, $false, $false, $false, $true, 1, $false, " a. numbering the fingers, establishing pattern for which finger will be landed on, creating a way for the number to be found.", 2) | Out-Null

# New paragraphs to append after that paragraph, in document order.
$newParagraphs = @(
    "",
    "4. ",
    " a. yes",
    " b. yes",
    "",
    "5. ",
    " a. number divide by 5 which is going to give you a number. If the number is greater than eight then subtract 8 until you get a number equal to or less than 8. For numbers 1-8",
    "1=5",
    "2=2",
    "3=3",
    "4=4",
    "5=1",
    "6=4",
    "7=3",
    "8=2",
    "For the numbers 1-5 that 1-8 equal, the numbers 1-5 are equal to fingers. ",
    "1=thumb",
    "2=index finger",
    "3=middle finger",
    "4=ring finger",
    "5=pinky finger",
    "So the answers to the questions are as follows",
    "a) 1-10 index finger",
    "b) 1-100 ring finger",
    "c) 1-1000 index finger"
)

foreach ($text in $newParagraphs) {
    $lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $lastPara.Range.InsertParagraphAfter() | Out-Null
    $newPara = $d.Paragraphs.Item($d.Paragraphs.Count)
    $newPara.Range.InsertXML((New-RunParagraphXml $text)) | Out-Null
}
